# The deck ships two embedded themes: theme1.xml ("Office Theme" colours,
# used by the Notes Master) and theme2.xml ("Integral" colours, used by the
# Slide Master / all slide layouts). The edit swaps which colour palette is
# applied where: the Slide Master now takes on the default Office colours
# while the Notes Master takes on the Integral colours.
#
# Drive it the way a user would from the Design tab: recolour the active
# theme's 12-slot colour scheme (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink)
# via ThemeColorScheme, swapping in the values that used to live in the
# other theme part.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# New values = the old "Office Theme" palette (previously theme1.xml),
# expressed as VBA RGB() longs (0x00BBGGRR).
$newColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = $newColors[$i - 1]
}
